# Weekly update: a new Choclo price record for "Región de Arica y Parinacota"
# is inserted as the new row 270 (pushing the existing rows 270-307 down to
# 271-308), and the sheet's used-range dimension grows from R307 to R308.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 270..307 down by one to make room for the new record.
$ws.Rows.Item(270).Insert()

# Populate the newly inserted row 270 with the new record.
$ws.Range("A270").Value = 9
$ws.Range("B270").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C270").Value = "Metropolitana"
$ws.Range("D270").Value = 44504
$ws.Range("E270").Value = 13
$ws.Range("F270").Value = 100112024
$ws.Range("G270").Value = "Choclo"
$ws.Range("H270").Value = "Dulce o Americano"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 36
$ws.Range("K270").Value = 39000
$ws.Range("L270").Value = 41000
$ws.Range("M270").Value = 40000
$ws.Range("N270").Value = "$/malla 70 unidades"
$ws.Range("O270").Value = "Región de Arica y Parinacota"
$ws.Range("P270").Value = 571
$ws.Range("Q270").Value = 70
$ws.Range("R270").Value = "Hortaliza"
